# parameters.xlsx -- "Changed lqr sim. Cleaned up runfile"
#
# 1) Re-enter the INDIRECT()/cross-reference formulas on the "Main" sheet as
#    single multi-cell Range.Formula assignments so Excel collapses each
#    column run (E2:E7, F2:F7, ... J2:J7, E9:E14, G9:G14, I9:I14) into one
#    shared formula.
# 2) Duplicate worksheet "63" to a new trailing worksheet "64" (sheetId 10),
#    tweak its three changed parameter values, and make it the active tab.

$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")

# --- Main!E2:J7 -----------------------------------------------------------
$main.Range("E2:E7").Formula = '=INDIRECT("''"&E$1&"''!A"&$D2)'
$main.Range("F2:F7").Formula = '=INDIRECT("''"&E$1&"''!B"&$D2)'
$main.Range("G2:G7").Formula = '=INDIRECT("''"&G$1&"''!A"&$D2)'
$main.Range("H2:H7").Formula = '=INDIRECT("''"&G$1&"''!B"&$D2)'
$main.Range("I2:I7").Formula = '=INDIRECT("''"&I$1&"''!A"&$D2)'
$main.Range("J2:J7").Formula = '=INDIRECT("''"&I$1&"''!B"&$D2)'

# --- Main!E9:J14 (mirrors the row-2..7 results) ----------------------------
$main.Range("E9:E14").Formula = '=F2'
$main.Range("G9:G14").Formula = '=H2'
$main.Range("I9:I14").Formula = '=J2'

# --- add sheet "64" as a copy of "63" --------------------------------------
$sheet63 = $wb.Worksheets.Item("63")
$sheet63.Copy([System.Reflection.Missing]::Value, $sheet63)

$sheet64 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet64.Name = "64"

$sheet64.Range("B1").Value = 64
$sheet64.Range("B10").Value = 0.45
$sheet64.Range("B23").Value = 100

$sheet64.Range("A1:B25").Select()
